# Commit: "fix: prevent hidden columns from being labeled upon detecting changes (#11)"
#
# The "Aenderung" (L) column flags rows where the FV2210 side (columns B:K)
# differs from the mirrored FV2304 side (columns M:V). A bug used to compare
# a column that should have been skipped (hidden), producing false-positive
# "changed" markers. This edit clears those false-positive markers (column L)
# for the affected rows, 120-197 (except 133 and 142, which have genuine
# differences and must keep their marker).
#
# Some of the cleared rows are also the first ("group header") row for a new
# record (the B column's value changes from the previous row): those rows
# additionally lose their distinct banding/highlight style across the whole
# A:V span, falling back to the plain "no-diff" group-header look used
# elsewhere in the sheet (same as e.g. row 2, 9, 14, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that are the first row of a new group: the whole A:V row is restyled
# to match the plain group-header appearance (copy format from row 2, which
# already has that look).
$fullChangeRows = @(120,123,128,131,135,140,144,147,151,155,159,163,171,173,177,181,189,192,195)

# Rows where only the L (Aenderung) marker/highlight needs to be removed.
$lOnlyRows = @(121,122,124,125,126,127,129,130,132,134,136,137,138,139,141,143,145,146,148,149,150,152,153,154,156,157,158,160,161,162,164,165,166,167,168,169,170,172,174,175,176,178,179,180,182,183,184,185,186,187,188,190,191,193,194,196,197)

$ws.Range("A2:V2").Copy()
foreach ($r in $fullChangeRows) {
    $ws.Range("A" + $r + ":V" + $r).PasteSpecial(-4122)  # xlPasteFormats
}

$ws.Range("L2").Copy()
foreach ($r in $lOnlyRows) {
    $ws.Range("L" + $r).PasteSpecial(-4122)  # xlPasteFormats
}

# Clear the stale "ÄNDERUNG" marker text now that the highlight is gone.
foreach ($r in ($fullChangeRows + $lOnlyRows)) {
    $ws.Range("L" + $r).Value = ""
}
